# Apply the "delete element from array" edit:
#   1. fix a stray upper-case "I" typo in the shift-left-on-delete pseudocode
#      (A[ I ] = A[ i-1 ];  ->  A[ i ] = A[ i-1 ];)
#   2. append a new "Eliminar" (array element deletion) section, with its own
#      pseudocode block and time-complexity note, after the existing content.
$d = $word.ActiveDocument

# --- Edit 1: "A[ I ] = A[ i-1 ];" -> "A[ i ] = A[ i-1 ];" -------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(" I ] = A[ i-1 ];", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find target text for edit 1"
}
$rng1.Expand(4) | Out-Null
$frag1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00E70F3A" w:rsidRDefault="00E70F3A" w:rsidP="001F0A20"><w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>A[</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> ] = A[ i-1 ];</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t>0 – n (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>depende</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>indice</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng1.InsertXML($frag1) | Out-Null

# --- Edit 2: split the "min(1)," run, then append the new section ----------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("La complejidad de tiempo es:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find target text for edit 2"
}
$rng2.Expand(4) | Out-Null
$frag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00E70F3A" w:rsidRPr="00E70F3A" w:rsidRDefault="00E70F3A" w:rsidP="001F0A20"><w:pPr><w:jc w:val="both"/></w:pPr><w:r w:rsidRPr="00E70F3A"><w:t>La complejidad de tiempo es:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>min(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">1), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>max</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(n) = O(n).</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Eliminar, </w:t></w:r><w:r><w:t xml:space="preserve">una vez eliminado el elemento, no se debe dejar el espacio </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>vacio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>asi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> que se debe mover todos los elementos posteriores hacia la izquierda. (el movimiento es hacia la izquierda hasta </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>length</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> -1)</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>X = A[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>index</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>];</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>1</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">For (int </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> = index; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> &lt;</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> length; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>++) {</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>A[</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> ] = A [ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> + 1];</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t>0-n</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Length--;</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t>1</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Complejidad de tiempo: min(2), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>max</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> =</w:t></w:r><w:r><w:t xml:space="preserve"> n + 2 = O(n).</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng2.InsertXML($frag2) | Out-Null

# InsertXML replacing the document's last paragraph leaves one extra blank
# paragraph behind (to hold the package's terminating paragraph mark); merge
# it back into the preceding paragraph, which is now the real last paragraph.
$lastPara = $d.Paragraphs.Last
$gap = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End)
$gap.Delete() | Out-Null

Write-Output "done"
